$d = $word.ActiveDocument

# 1) Merge split runs: "Worked with purchasing ... suppliers" + "while giving technical direction ... needs."
$old1 = "Worked with purchasing and technical leaders to evaluate virtual ECU platforms and suppliers while giving technical direction on integration and model development needs."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# 2) Merge split runs: "Employed technical knowledge of VECUs ... customer " + "needs" + "."
$old2 = "Employed technical knowledge of VECUs to efficiently troubleshoot and resolve open issues, prioritizing customer needs."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# 3) Merge split runs: "Leveraged expertise in " + "virtual" + " workflows to collaborate ... value stream."
$old3 = "Leveraged expertise in virtual workflows to collaborate with product owners and release train engineers, ensuring successful integration into the existing SAFe value stream."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null

# 4) OTA bullet: insert "Chevy " before "Volts and Bolts)."
$old4 = "OTA (“Over-the-Air”) Release Engineer for existing OnStar connected customers for PHEV and Evs (Volts and Bolts)."
$new4 = "OTA (“Over-the-Air”) Release Engineer for existing OnStar connected customers for PHEV and Evs (Chevy Volts and Bolts)."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# 5) Agile bullet: append clause about giving presentations to upper management.
$old5 = "Employed Agile and Scrum principles to preempt and quickly resolve issues."
$new5 = "Employed Agile and Scrum principles to preempt and quickly resolve issues while giving presentations to upper management for decisions."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

# 6) Merge split runs: "Automation solution developer supporting " + "electronic module" + " production, and launch activities for manufacturing and testing services group."
$old6 = "Automation solution developer supporting electronic module production, and launch activities for manufacturing and testing services group."
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $old6, 2) | Out-Null

# 7) Merge split runs: "Designed and documented manufacturing " + "test solutions using UML methods and engaged in peer-review to verify " + "functionality."
$old7 = "Designed and documented manufacturing test solutions using UML methods and engaged in peer-review to verify functionality."
$d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $old7, 2) | Out-Null

# 8) Merge split runs: "Led " + "identif" + "ication" + " and mitigat" + "ion of" + " risks in production lines " + "during" + " PFMEA discussions with a global cross-functional team."
$old8 = "Led identification and mitigation of risks in production lines during PFMEA discussions with a global cross-functional team."
$d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $old8, 2) | Out-Null

# 9) Merge split runs: "Collaborated with design ... related to " + "electronics modules" + "."
$old9 = "Collaborated with design and quality engineers to develop clear work instructions and troubleshooting procedures for resolving production problems related to electronics modules."
$d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $old9, 2) | Out-Null

# 10) Merge split runs: "Drove the integration " + "and " + "extension of " + "HMI" + " touchscreen automation tool ... testing time."
$old10 = "Drove the integration and extension of HMI touchscreen automation tool created in Python. This enabled our lab to perform automated testing at scale for features with Ford Sync screen for the first time: saving 100s of hours of human testing time."
$d.Content.Find.Execute($old10, $true, $false, $false, $false, $false, $true, 1, $false, $old10, 2) | Out-Null

# 11) Insert new bullet after "Led BOM validation..." reading "Gave automation solution demos to management."
$old11 = "Led BOM validation and parts acquisition for HIL Simulator subsystem benches concurrently for multiple programs."
$rng11 = $d.Content
$found11 = $rng11.Find.Execute($old11, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found11) {
    $para11 = $rng11.Paragraphs(1)
    $para11.Range.InsertParagraphAfter()
    $newPara11 = $para11.Next()
    $newPara11.Range.InsertAfter("Gave automation solution demos to management.")
}

# 12) Flip w:overflowPunct false -> true on the Normal, TOC Heading, and No Spacing styles.
$styles = $d.Styles
$overflowPunctStyleNames = @("Normal", "TOC Heading", "No Spacing")
foreach ($styleName in $overflowPunctStyleNames) {
    $style = $styles.Item($styleName)
    $style.ParagraphFormat.HangingPunctuation = $true
}
